$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Hoja1")

# The account-statement table (rows 16-22) lists the most recent 7 "Periodo Mora"
# periods together with their "Valor Mora" / "Salario Basico" figures. The
# workbook's feeder macro dropped the oldest period (2404) and rolled in the
# newest one (2410), so the period column now reads newest-to-oldest and the
# "Valor Mora" figure that used to belong to period 2404 now travels with it
# to the bottom row, while period 2410's smaller first/partial-month value
# moves up to the top row.

$ws.Range("E16").Value = "2410"
$ws.Range("E17").Value = "2409"
$ws.Range("E18").Value = "2408"
$ws.Range("E19").Value = "2407"
$ws.Range("E20").Value = "2406"
$ws.Range("E21").Value = "2405"
$ws.Range("E22").Value = "2404"

$ws.Range("F16").Value = 54600
$ws.Range("F22").Value = 91000
